# Project number / comments correction in the PDR article index.
#
# D2: replace the short MOG blurb with the expanded write-up (priors,
#     regularization, loss function breakdown).
# E2: add the "relevant youtube video" link that was missing before.
# Selection moves from C2 to E2 to reflect the newly-edited cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "The article describes a MOG model with spatial priors (probabilities are vary over voxels). It also takes into considoration distoration and regularization (optimization of I voxals with k Geussians each).`nSegmenting magnetic resonance (MR) images into different tissue classes, using a modified Gaussian mixture model.`nMOG - Mixture of Gaussians distribution model. This is parametric representations of image intensity distributions.`nLoss function conssists of MOG parameters (meo,sigma,gama) , deforamtion (infectes the spatial priors - alpha), bias corrections (beta). Optimization using: EM for MOG parameters. LM for deforamation and bias.`n"

$ws.Range("E2").Value = "a relevant youtube video list:`nhttps://www.youtube.com/watch?v=REypj2sy_5U&index=1&list=PLBv09BD7ez_4e9LtmK626Evn1ion6ynrt"

$ws.Range("E2").Select()
